# Insert a new weekly data row for "Poroto verde" at row 36, pushing the
# existing rows 36:91 down to 37:92 (matches the commit's "Fruta / hortaliza,
# semanal" weekly-refresh pattern: a new week's record is prepended and the
# older history shifts down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 36 through 91 down by one row.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with this week's record.
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 45162
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112031
$ws.Range("G36").Value = "Poroto verde"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 1800
$ws.Range("K36").Value = 1100
$ws.Range("L36").Value = 1200
$ws.Range("M36").Value = 1144
$ws.Range("N36").Value = "`$/kilo"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 1144
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"
